# Update "paises.xlsx" (countries COVID leaderboard) per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last refreshed" timestamp banner.
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 12:05"

# Bielorrusia overtakes Ecuador in total cases (rows 25/26 swap, with
# Bielorrusia getting fresh numbers and Ecuador keeping its former stats).
$ws.Range("A25").Value = "Bielorrusia"
$ws.Range("B25").Value = 38059
$ws.Range("C25").Value = 915
$ws.Range("D25").Value = 15086
$ws.Range("E25").Value = 22765
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 208

$ws.Range("A26").Value = "Ecuador"
$ws.Range("B26").Value = 37355
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 18003
$ws.Range("E26").Value = 16149
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 3203

# Emiratos Arabes Unidos overtakes Portugal and Suiza (rows 30/31/32 shift
# down, Emiratos gets fresh numbers).
$ws.Range("A30").Value = "Emiratos Arabes Unidos"
$ws.Range("B30").Value = 31086
$ws.Range("C30").Value = 779
$ws.Range("D30").Value = 15982
$ws.Range("E30").Value = 14851
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 253

$ws.Range("A31").Value = "Portugal"
$ws.Range("B31").Value = 30788
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 17822
$ws.Range("E31").Value = 11636
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 1330

$ws.Range("A32").Value = "Suiza"
$ws.Range("B32").Value = 30746
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 28200
$ws.Range("E32").Value = 633
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 1913

# Kuwait overtakes Colombia (rows 36/37 swap, Kuwait gets fresh numbers).
$ws.Range("A36").Value = "Kuwait"
$ws.Range("B36").Value = 22575
$ws.Range("C36").Value = 608
$ws.Range("D36").Value = 7306
$ws.Range("E36").Value = 15097
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 7
$ws.Range("H36").Value = 172

$ws.Range("A37").Value = "Colombia"
$ws.Range("B37").Value = 21981
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 5265
$ws.Range("E37").Value = 15966
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 750

# Rumania's own stats refresh (no reordering here).
$ws.Range("B40").Value = 18429
$ws.Range("C40").Value = 146
$ws.Range("D40").Value = 11874
$ws.Range("E40").Value = 5345
